$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2025-07-11T12:29:53+00:00"

# Set Jurisdiction value (row 11, column B) to FRANCE
$ws.Range("B11").Value = "FRANCE"
